$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("4MC")

# Label the new summary rows
$ws.Range("G38").Value = "Minimum Distance"
$ws.Range("G39").Value = "Assigned Cluster"
$ws.Range("G39").Font.Bold = $true

# Row 38: minimum of the four "Distance to Cluster" rows (34:37) for each campaign column
$ws.Range("L38:DG38").Formula = "=MIN(L34:L37)"

# Row 39: which cluster (1-4) achieved that minimum distance
$ws.Range("L39:DG39").Formula = "=MATCH(L38,L34:L37,0)"

$ws.Activate() | Out-Null
$ws.Range("L39:DG39").Select() | Out-Null
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 107
